$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are numeric-looking text (prices like "1.001", "28.042.60")
# that must remain stored as text, matching the original inlineStr cells.
# Force text format, assign value, then restore default style so no stray
# cell style is introduced.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '28.042.60'
$ws.Range('E2').Value = '  -1.29%  '
Set-TextValue 'D3' '1.793.26'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  -0.01%  '
Set-TextValue 'D5' '316.80'
$ws.Range('E5').Value = '  +0.80%  '
Set-TextValue 'D6' '1.001'
$ws.Range('E6').Value = '  +0.01%  '
Set-TextValue 'D7' '0.5358'
$ws.Range('E7').Value = '  -1.79%  '
Set-TextValue 'D8' '0.3766'
$ws.Range('E8').Value = '  -1.50%  '
Set-TextValue 'D9' '0.07428'
$ws.Range('E9').Value = '  -1.87%  '
Set-TextValue 'D10' '41.92'
$ws.Range('E10').Value = '  -2.00%  '
Set-TextValue 'D11' '1.093'
$ws.Range('E11').Value = '  -2.82%  '
Set-TextValue 'D12' '1.000'
$ws.Range('E12').Value = '  -0.01%  '
Set-TextValue 'D13' '20.56'
$ws.Range('E13').Value = '  -2.74%  '
Set-TextValue 'D14' '6.121'
$ws.Range('E14').Value = '  -1.08%  '
Set-TextValue 'D15' '7.238'
$ws.Range('E15').Value = '  -2.12%  '
Set-TextValue 'D16' '1.787.42'
$ws.Range('E16').Value = '  -0.76%  '
Set-TextValue 'D17' '89.06'
$ws.Range('E17').Value = '  -2.71%  '
$ws.Range('E18').Value = '  -1.29%  '
Set-TextValue 'D19' '0.06497'
$ws.Range('E19').Value = '  +0.74%  '
$ws.Range('E20').Value = '  -0.02%  '
Set-TextValue 'D21' '17.26'
$ws.Range('E21').Value = '  -0.38%  '
Set-TextValue 'D22' '5.894'
$ws.Range('E22').Value = '  -1.04%  '
Set-TextValue 'D23' '28.064.74'
$ws.Range('E23').Value = '  -1.29%  '
Set-TextValue 'D24' '11.16'
$ws.Range('E24').Value = '  -2.70%  '
Set-TextValue 'D25' '2.090'
$ws.Range('E25').Value = '  -1.87%  '
Set-TextValue 'D26' '155.35'
$ws.Range('E26').Value = '  -2.58%  '
Set-TextValue 'D27' '20.30'
$ws.Range('E27').Value = '  -1.99%  '
Set-TextValue 'D28' '1.993.25'
$ws.Range('E28').Value = '  -0.68%  '
Set-TextValue 'D29' '2.307'
$ws.Range('E29').Value = '  -3.66%  '
Set-TextValue 'D30' '121.14'
$ws.Range('E30').Value = '  -1.90%  '
Set-TextValue 'D31' '1.119'
$ws.Range('E31').Value = '  -0.64%  '
Set-TextValue 'D32' '0.1062'
$ws.Range('E32').Value = '  +3.82%  '
$ws.Range('E33').Value = '  -0.38%  '
Set-TextValue 'D34' '5.558'
$ws.Range('E34').Value = '  -3.45%  '
$ws.Range('E35').Value = '  -4.34%  '
Set-TextValue 'D36' '0.06487'
$ws.Range('E36').Value = '  -4.59%  '
Set-TextValue 'D37' '0.02289'
$ws.Range('E37').Value = '  -1.54%  '
Set-TextValue 'D38' '5.023'
$ws.Range('E38').Value = '  -2.57%  '
Set-TextValue 'D39' '8.460'
$ws.Range('E39').Value = '  -3.54%  '
Set-TextValue 'D40' '0.6185'
$ws.Range('E40').Value = '  -3.32%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D41' '11.15'
$ws.Range('E41').Value = '  -4.42%  '
$ws.Range('B42').Value = 'WEMIXTOKEN'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D42' '1.449'
$ws.Range('E42').Value = '  +1.79%  '
$ws.Range('E43').Value = '  +1.94%  '
Set-TextValue 'D44' '13.37'
$ws.Range('E44').Value = '  -1.92%  '
$ws.Range('E45').Value = '  +0.06%  '
Set-TextValue 'D46' '0.5785'
$ws.Range('E46').Value = '  -3.31%  '
Set-TextValue 'D47' '124.90'
$ws.Range('E47').Value = '  -1.22%  '
Set-TextValue 'D48' '1.189'
$ws.Range('E48').Value = '  +3.31%  '
Set-TextValue 'D49' '1.927'
$ws.Range('E49').Value = '  -3.88%  '
Set-TextValue 'D50' '0.06821'
$ws.Range('E50').Value = '  -2.09%  '
Set-TextValue 'D51' '71.29'
$ws.Range('E51').Value = '  -2.29%  '
